$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6264044943820225
$ws1.Range("C2").Value = 0.5816203143893591
$ws1.Range("D2").Value = 0.900749063670412
$ws1.Range("E2").Value = 0.7068332108743571
$ws1.Range("F2").Value = 0.8116773540330746
$ws1.Range("G2").Value = 0.8821330323763843
$ws1.Range("H2").Value = 0.7688107562176492
$ws1.Range("I2").Value = 481
$ws1.Range("J2").Value = 346
$ws1.Range("K2").Value = 188
$ws1.Range("L2").Value = 53

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.7800829875518672
$ws2.Range("C2").Value = 0.352059925093633
$ws2.Range("D2").Value = 0.4851612903225807

$ws2.Range("B3").Value = 0.5816203143893591
$ws2.Range("C3").Value = 0.900749063670412
$ws2.Range("D3").Value = 0.7068332108743571

$ws2.Range("B4").Value = 0.6264044943820225
$ws2.Range("C4").Value = 0.6264044943820225
$ws2.Range("D4").Value = 0.6264044943820225
$ws2.Range("E4").Value = 0.6264044943820225

$ws2.Range("B5").Value = 0.6808516509706132
$ws2.Range("C5").Value = 0.6264044943820225
$ws2.Range("D5").Value = 0.5959972505984689

$ws2.Range("B6").Value = 0.6808516509706132
$ws2.Range("C6").Value = 0.6264044943820225
$ws2.Range("D6").Value = 0.5959972505984689

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 188
$ws3.Range("C2").Value = 346
$ws3.Range("B3").Value = 53
$ws3.Range("C3").Value = 481
